$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting for column A (date column) from the last existing row
$ws.Range("A343").Copy()
$ws.Range("A344:A357").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$rows = @(
  @(344, 44418, 1, 19, 57.60543310190099),
  @(345, 44419, 0, 19, 57.60543310190099),
  @(346, 44420, 4, 13, 39.41424370130067),
  @(347, 44421, 1, 14, 42.44610860140072),
  @(348, 44422, 3, 12, 36.38237880120062),
  @(349, 44423, 1, 13, 39.41424370130067),
  @(350, 44424, 6, 16, 48.50983840160082),
  @(351, 44425, 1, 16, 48.50983840160082),
  @(352, 44426, 1, 17, 51.54170330170087),
  @(353, 44427, 0, 13, 39.41424370130067),
  @(354, 44428, 10, 22, 66.70102780220114),
  @(355, 44429, 3, 22, 66.70102780220114),
  @(356, 44430, 2, 23, 69.73289270230119),
  @(357, 44431, 7, 24, 72.76475760240123)
)

foreach ($r in $rows) {
  $rowNum = $r[0]
  $ws.Cells.Item($rowNum, 1).Value = $r[1]
  $ws.Cells.Item($rowNum, 2).Value = $r[2]
  $ws.Cells.Item($rowNum, 3).Value = $r[3]
  $ws.Cells.Item($rowNum, 4).Value = $r[4]
}
